# Update "想去人数" (F column) counts and a few status changes on both the
# "展览" sheet (index 1) and the "全部类型" sheet (index 4), which mirror the
# same data set in this workbook.

$wb = $excel.ActiveWorkbook

$targetSheets = @(1, 4)

foreach ($sheetIndex in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # Simple numeric updates to column F (想去人数 / "want to go" count)
    $ws.Cells.Item(2, 6).Value  = 1842   # F2:  1833 -> 1842
    $ws.Cells.Item(7, 6).Value  = 1520   # F7:  1510 -> 1520
    $ws.Cells.Item(8, 6).Value  = 19     # F8:  18   -> 19
    $ws.Cells.Item(9, 6).Value  = 595    # F9:  593  -> 595
    $ws.Cells.Item(10, 6).Value = 363    # F10: 362  -> 363
    $ws.Cells.Item(14, 6).Value = 220    # F14: 219  -> 220
    $ws.Cells.Item(16, 6).Value = 133    # F16: 132  -> 133
    $ws.Cells.Item(17, 6).Value = 103    # F17: 102  -> 103
    $ws.Cells.Item(19, 6).Value = 3575   # F19: 3563 -> 3575
    $ws.Cells.Item(20, 6).Value = 426    # F20: 425  -> 426
    $ws.Cells.Item(21, 6).Value = 320    # F21: 318  -> 320
    $ws.Cells.Item(23, 6).Value = 128    # F23: 120  -> 128
    $ws.Cells.Item(26, 6).Value = 1339   # F26: 1320 -> 1339
    $ws.Cells.Item(27, 6).Value = 136    # F27: 135  -> 136

    # Row 11: event got cancelled -> name suffix + ticket price becomes
    # "not for sale" (string) instead of a numeric price.
    $ws.Cells.Item(11, 3).Value = "南昌·童话镇国乙&鸢only（取消）"
    $ws.Cells.Item(11, 7).Value = "不可售"
}
